$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LIST")

# Row 2: A2 changes value (shared string content change)
$ws.Range("A2").Value = "AD.SEC.001.FON.01"

# Row 3: A3 changes value and loses its style (reset to default/General style)
$ws.Range("A3").Value = "RO.ACT.003HAB.SRM"
$ws.Range("A3").Style = "Normal"

# Row 3: B3 gets a new numeric value
$ws.Range("B3").Value = 20

# Row 4: A4 cleared but keeps its style/format
$ws.Range("A4").ClearContents()

# Row 5: A5 cleared but keeps its style/format
$ws.Range("A5").ClearContents()

# Row 6: A6 cell removed entirely
$ws.Range("A6").Clear()

# Update selection to B4
$ws.Range("B4").Select()
